# Update the pharmacy report: remove the "RICHI PANTHENOL CREAM 50GM" line item
# (row 10) from the product table, which shifts the following product rows up
# by one, and refresh the dependent totals / footer accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire worksheet row holding "RICHI PANTHENOL CREAM 50GM".
# This shifts rows 11-15 up to rows 10-14, re-flows the merged cell ranges,
# and drops the now-unused shared string automatically on save.
$ws.Rows.Item(10).Delete()

# The sequential item-number column ("م") must keep reading 1..9 for the
# remaining products, so restore the numbering for the rows that moved up.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(12, 1).Value = 9

# The grand-total cell (column K of the totals row, now row 13) is a static
# value, not a formula, so subtract the removed item's amount (62) by hand:
# 376.36 - 62 = 314.36
$ws.Cells.Item(13, 11).Value = 314.36

# The footer row (timestamp / page / developed-by) moved from row 15 to row
# 14 and regains its original (slightly taller) row height.
$ws.Rows.Item(14).RowHeight = 17.25
